$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new data row (row 33), mirroring the pattern of the existing rows.
$ws.Range("A33").Value = 10002
$ws.Range("B33").Value = 10032
$ws.Range("C33").Value = "eng"
$ws.Range("D33").Value = $true
$ws.Range("E33").Value = "superadmin"
$ws.Range("F33").Value = "now()"
$ws.Range("G33").Value = "now()"

# Update the active selection to match the edited workbook state.
$ws.Range("B30").Select()
